# Refresh the cryptos price/volume table with the latest scraped values
# (GitHub Actions cron update). Cells are plain text in the source data
# (e.g. "26.906.71", "  +0.82%  "), so for Price-column values that look
# like plain numbers we briefly force Text format before assigning, then
# restore General, to avoid Excel silently coercing them into numeric
# values (which would lose formatting such as trailing zeros / truncate
# "19.84" style figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.906.71'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.643.75'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.72%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.94'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('E8').Value = '  +1.79%  '
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.84'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('E10').Value = '  +4.67%  '
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.13'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.600.10'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.40'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('D17').Value = '26.901.99'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '220.06'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('E19').Value = '  +4.29%  '
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.64'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('E21').Value = '  +7.47%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.39'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('E23').Value = '  +4.68%  '
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.09'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.42'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('E27').Value = '  +6.01%  '
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.83'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.248.01'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.43'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0175'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.539'
$ws.Range('D38').NumberFormat = "General"
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.836'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('E39').Value = '  +4.22%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.805'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D43').Value = '1.783.51'
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('E44').Value = '  -2.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '60.84'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.55'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +18.33%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0514'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0974'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.60'
$ws.Range('D51').NumberFormat = "General"
$ws.Range('E51').Value = '  +2.49%  '
